$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Text = ""
